$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: endpoint path now namespaced under /api, description added
$ws.Range("A18").Value = "/api/ngo/:email"
$ws.Range("B18").Value = "delete"
$ws.Range("C18").Value = "delete the ngo"

# Row 19: endpoint path now namespaced under /api, description added
$ws.Range("A19").Value = "/api/user/:email"
$ws.Range("B19").Value = "delete"
$ws.Range("C19").Value = "delete the user"

# Row 20 (new): list unverified users endpoint
$ws.Range("A20").Value = "/api/unverifiedusers"
$ws.Range("B20").Value = "get"
$ws.Range("C20").Value = "list of unverified users"

$ws.Range("A20").Style = $ws.Range("A19").Style

$ws.Range("D21").Select()
$excel.ActiveWindow.ScrollRow = 7
